$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 10 (shifts old summary row 10 -> row 12,
# preserving its original cell types, e.g. the blank-text B10:D10 cells).
$ws.Rows("10:11").Insert()

# Restore the number-style (border/bold/centered) on the new A10:A11 index cells
# to match the rest of column A (copy direct formatting from A9).
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update changed values in the existing data rows (2-9) ---
$ws.Cells.Item(2, 5).Value = 0.75
$ws.Cells.Item(2, 6).Value = 0.5
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 0.8181818181818182
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 0.5
$ws.Cells.Item(2, 11).Value = 0.8181818181818182
$ws.Cells.Item(2, 12).Value = 0.8
$ws.Cells.Item(2, 14).Value = 0.8
$ws.Cells.Item(2, 15).Value = 0.8775510204081632
$ws.Cells.Item(2, 16).Value = 0.869281045751634
$ws.Cells.Item(2, 17).Value = 0.8888888888888888
$ws.Cells.Item(2, 18).Value = 0.6
$ws.Cells.Item(2, 19).Value = 0.6666666666666666
$ws.Cells.Item(2, 21).Value = 0.3333333333333333
$ws.Cells.Item(2, 22).Value = 1
$ws.Cells.Item(2, 23).Value = 1
$ws.Cells.Item(2, 24).Value = 1

$ws.Cells.Item(3, 5).Value = 0.5
$ws.Cells.Item(3, 7).Value = 0.25
$ws.Cells.Item(3, 8).Value = 0.8181818181818182
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 11).Value = 0.9090909090909091
$ws.Cells.Item(3, 12).Value = 0.8
$ws.Cells.Item(3, 13).Value = 0.6666666666666666
$ws.Cells.Item(3, 14).Value = 0.8666666666666667
$ws.Cells.Item(3, 15).Value = 0.8163265306122449
$ws.Cells.Item(3, 16).Value = 0.8300653594771242
$ws.Cells.Item(3, 17).Value = 0.7777777777777778
$ws.Cells.Item(3, 18).Value = 0.4
$ws.Cells.Item(3, 19).Value = 0.3333333333333333
$ws.Cells.Item(3, 21).Value = 0.3333333333333333
$ws.Cells.Item(3, 22).Value = 1
$ws.Cells.Item(3, 23).Value = 0.8571428571428571
$ws.Cells.Item(3, 24).Value = 0.75

$ws.Cells.Item(4, 5).Value = 0.5
$ws.Cells.Item(4, 6).Value = 0.5
$ws.Cells.Item(4, 7).Value = 0.75
$ws.Cells.Item(4, 8).Value = 0.7272727272727273
$ws.Cells.Item(4, 9).Value = 0.8
$ws.Cells.Item(4, 10).Value = 0.5
$ws.Cells.Item(4, 11).Value = 0.8181818181818182
$ws.Cells.Item(4, 12).Value = 0.6
$ws.Cells.Item(4, 14).Value = 0.5333333333333333
$ws.Cells.Item(4, 15).Value = 0.7755102040816326
$ws.Cells.Item(4, 16).Value = 0.8104575163398693
$ws.Cells.Item(4, 17).Value = 0.7777777777777778
$ws.Cells.Item(4, 18).Value = 0.6
$ws.Cells.Item(4, 19).Value = 0.6666666666666666
$ws.Cells.Item(4, 21).Value = 0.3333333333333333
$ws.Cells.Item(4, 22).Value = 1
$ws.Cells.Item(4, 23).Value = 0.7142857142857143
$ws.Cells.Item(4, 24).Value = 1

$ws.Cells.Item(5, 2).Value = "hisditonly"
$ws.Cells.Item(5, 3).Value = "raw"
$ws.Cells.Item(5, 5).Value = 0.25
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 0.25
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 9).Value = 0.8
$ws.Cells.Item(5, 11).Value = 0.9090909090909091
$ws.Cells.Item(5, 12).Value = 0.6
$ws.Cells.Item(5, 13).Value = 0.6666666666666666
$ws.Cells.Item(5, 14).Value = 0.8
$ws.Cells.Item(5, 15).Value = 0.8775510204081632
$ws.Cells.Item(5, 16).Value = 0.7973856209150327
$ws.Cells.Item(5, 17).Value = 0.7777777777777778
$ws.Cells.Item(5, 18).Value = 0.4
$ws.Cells.Item(5, 19).Value = 0.6666666666666666
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 0.3333333333333333
$ws.Cells.Item(5, 22).Value = 1
$ws.Cells.Item(5, 23).Value = 0.5714285714285714

$ws.Cells.Item(6, 2).Value = "flowbot"
$ws.Cells.Item(6, 5).Value = 0.75
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = 1
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 0.9090909090909091
$ws.Cells.Item(6, 12).Value = 0.8
$ws.Cells.Item(6, 13).Value = 0.3333333333333333
$ws.Cells.Item(6, 14).Value = 1
$ws.Cells.Item(6, 15).Value = 0.9591836734693876
$ws.Cells.Item(6, 16).Value = 0.915032679738562
$ws.Cells.Item(6, 17).Value = 0.8888888888888888
$ws.Cells.Item(6, 18).Value = 1
$ws.Cells.Item(6, 19).Value = 0.6666666666666666
$ws.Cells.Item(6, 20).Value = 1
$ws.Cells.Item(6, 22).Value = 1
$ws.Cells.Item(6, 23).Value = 1
$ws.Cells.Item(6, 24).Value = 1

$ws.Cells.Item(7, 2).Value = "dit"
$ws.Cells.Item(7, 5).Value = 0.5
$ws.Cells.Item(7, 7).Value = 0.5
$ws.Cells.Item(7, 8).Value = 1
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 1
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.8
$ws.Cells.Item(7, 13).Value = 0.3333333333333333
$ws.Cells.Item(7, 14).Value = 0.9333333333333332
$ws.Cells.Item(7, 15).Value = 0.9183673469387755
$ws.Cells.Item(7, 16).Value = 0.869281045751634
$ws.Cells.Item(7, 17).Value = 0.7777777777777778
$ws.Cells.Item(7, 18).Value = 0.8
$ws.Cells.Item(7, 19).Value = 0.3333333333333333
$ws.Cells.Item(7, 21).Value = 0.6666666666666666
$ws.Cells.Item(7, 22).Value = 1
$ws.Cells.Item(7, 23).Value = 0.7142857142857143
$ws.Cells.Item(7, 24).Value = 1

$ws.Cells.Item(8, 2).Value = "pndit"
$ws.Cells.Item(8, 5).Value = 0.5
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 0.75
$ws.Cells.Item(8, 8).Value = 0.7272727272727273
$ws.Cells.Item(8, 9).Value = 1
$ws.Cells.Item(8, 10).Value = 1
$ws.Cells.Item(8, 11).Value = 0.9090909090909091
$ws.Cells.Item(8, 12).Value = 0.8
$ws.Cells.Item(8, 13).Value = 0.6666666666666666
$ws.Cells.Item(8, 14).Value = 0.6
$ws.Cells.Item(8, 15).Value = 0.9183673469387756
$ws.Cells.Item(8, 16).Value = 0.8823529411764706
$ws.Cells.Item(8, 17).Value = 0.8888888888888888
$ws.Cells.Item(8, 18).Value = 0.8
$ws.Cells.Item(8, 19).Value = 0.3333333333333333
$ws.Cells.Item(8, 22).Value = 1
$ws.Cells.Item(8, 23).Value = 1
$ws.Cells.Item(8, 24).Value = 1

$ws.Cells.Item(9, 2).Value = "hisdit"
$ws.Cells.Item(9, 7).Value = 0.75
$ws.Cells.Item(9, 8).Value = 0.9090909090909091
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.6
$ws.Cells.Item(9, 15).Value = 0.9183673469387756
$ws.Cells.Item(9, 16).Value = 0.954248366013072
$ws.Cells.Item(9, 17).Value = 0.7777777777777778
$ws.Cells.Item(9, 18).Value = 0.8

# --- Populate the two newly-inserted rows (10-11) ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "pnhisdit"
$ws.Cells.Item(10, 3).Value = "sgp"
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0.75
$ws.Cells.Item(10, 6).Value = 0.5
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 0.7272727272727273
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = 1
$ws.Cells.Item(10, 11).Value = 0.9090909090909091
$ws.Cells.Item(10, 12).Value = 0.8
$ws.Cells.Item(10, 13).Value = 0.3333333333333333
$ws.Cells.Item(10, 14).Value = 1
$ws.Cells.Item(10, 15).Value = 0.9795918367346941
$ws.Cells.Item(10, 16).Value = 0.9673202614379085
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = 1
$ws.Cells.Item(10, 19).Value = 1
$ws.Cells.Item(10, 20).Value = 0
$ws.Cells.Item(10, 21).Value = 0.6666666666666666
$ws.Cells.Item(10, 22).Value = 1
$ws.Cells.Item(10, 23).Value = 0.8571428571428571
$ws.Cells.Item(10, 24).Value = 1

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "hisditonly"
$ws.Cells.Item(11, 3).Value = "sgp"
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0.75
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0.25
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(11, 10).Value = 1
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.8
$ws.Cells.Item(11, 13).Value = 0.6666666666666666
$ws.Cells.Item(11, 14).Value = 0.7333333333333333
$ws.Cells.Item(11, 15).Value = 0.9183673469387756
$ws.Cells.Item(11, 16).Value = 0.8823529411764706
$ws.Cells.Item(11, 17).Value = 0.8888888888888888
$ws.Cells.Item(11, 18).Value = 0.8
$ws.Cells.Item(11, 19).Value = 1
$ws.Cells.Item(11, 20).Value = 0
$ws.Cells.Item(11, 21).Value = 0.3333333333333333
$ws.Cells.Item(11, 22).Value = 1
$ws.Cells.Item(11, 23).Value = 0.5714285714285714
$ws.Cells.Item(11, 24).Value = 0.25

